$wb = $excel.ActiveWorkbook

# --- Sheet "area_mixre": update descriptive stats for the "area" column ---
$wsMixre = $wb.Worksheets.Item("area_mixre")
$wsMixre.Range("B2").Value = 190
$wsMixre.Range("B3").Value = 4.112423558429414
$wsMixre.Range("B4").Value = 3.804264655101014
$wsMixre.Range("B5").Value = 0.1718452795435217
$wsMixre.Range("B6").Value = 1.269746006822269
$wsMixre.Range("B7").Value = 3.034135648755669
$wsMixre.Range("B8").Value = 5.919150946757894

# --- Sheet "area_pop_sum": update population and density totals ---
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")
$wsPopSum.Range("B3").Value = 902158
$wsPopSum.Range("B4").Value = 1154.598919330731
